$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 40.013281
$ws.Cells.Item(2, 8).Value = 120.039843
$ws.Cells.Item(2, 9).Value = 0.9259567983352626
$ws.Cells.Item(2, 10).Value = 0.9352019795456475
$ws.Cells.Item(2, 13).Value = 36.48539666666667
$ws.Cells.Item(2, 14).Value = 109.45619
$ws.Cells.Item(2, 15).Value = 0.4260639713374229
$ws.Cells.Item(2, 16).Value = 0.4324607845540777
$ws.Cells.Item(2, 17).Value = 1459.900429219797
$ws.Cells.Item(2, 18).Value = 13139.10386297817
$ws.Cells.Item(2, 19).Value = 0.3945168307856071
$ws.Cells.Item(2, 20).Value = 0.4044381817908372
$ws.Cells.Item(3, 7).Value = 40.013281
$ws.Cells.Item(3, 8).Value = 120.039843
$ws.Cells.Item(3, 9).Value = 0.9259567983352626
$ws.Cells.Item(3, 10).Value = 0.9352019795456475
$ws.Cells.Item(3, 15).Value = 0.1743777127077069
$ws.Cells.Item(3, 16).Value = 0.1769957741547643
$ws.Cells.Item(3, 17).Value = 597.5020531053934
$ws.Cells.Item(3, 18).Value = 5377.518477948541
$ws.Cells.Item(3, 19).Value = 0.1614662285598545
$ws.Cells.Item(3, 20).Value = 0.1655267983607499
$ws.Cells.Item(4, 7).Value = 40.013281
$ws.Cells.Item(4, 8).Value = 120.039843
$ws.Cells.Item(4, 9).Value = 0.9259567983352626
$ws.Cells.Item(4, 10).Value = 0.9352019795456475
$ws.Cells.Item(4, 13).Value = 9.680823666666667
$ws.Cells.Item(4, 14).Value = 29.042471
$ws.Cells.Item(4, 15).Value = 0.1130493445068016
$ws.Cells.Item(4, 16).Value = 0.1147466378470605
$ws.Cells.Item(4, 17).Value = 387.3615176857837
$ws.Cells.Item(4, 18).Value = 3486.253659172053
$ws.Cells.Item(4, 19).Value = 0.1046788090934181
$ws.Cells.Item(4, 20).Value = 0.1073112828607785
$ws.Cells.Item(5, 7).Value = 40.013281
$ws.Cells.Item(5, 8).Value = 120.039843
$ws.Cells.Item(5, 9).Value = 0.9259567983352626
$ws.Cells.Item(5, 10).Value = 0.9352019795456475
$ws.Cells.Item(5, 13).Value = 3.79999
$ws.Cells.Item(5, 14).Value = 7.59998
$ws.Cells.Item(5, 15).Value = 0.04437498227672168
$ws.Cells.Item(5, 16).Value = 0.0300274777826206
$ws.Cells.Item(5, 17).Value = 152.05006766719
$ws.Cells.Item(5, 18).Value = 912.3004060031401
$ws.Cells.Item(5, 19).Value = 0.04108931651513723
$ws.Cells.Item(5, 20).Value = 0.02808175666306974
$ws.Cells.Item(6, 7).Value = 40.013281
$ws.Cells.Item(6, 8).Value = 120.039843
$ws.Cells.Item(6, 9).Value = 0.9259567983352626
$ws.Cells.Item(6, 10).Value = 0.9352019795456475
$ws.Cells.Item(6, 13).Value = 20.734808
$ws.Cells.Item(6, 14).Value = 62.204424
$ws.Cells.Item(6, 15).Value = 0.242133989171347
$ws.Cells.Item(6, 16).Value = 0.245769325661477
$ws.Cells.Item(6, 17).Value = 829.667698985048
$ws.Cells.Item(6, 18).Value = 7467.009290865432
$ws.Cells.Item(6, 19).Value = 0.2242056133812456
$ws.Cells.Item(6, 20).Value = 0.2298439598702122
$ws.Cells.Item(7, 9).Value = 0.03933964692088724
$ws.Cells.Item(7, 10).Value = 0.0397324321622614
$ws.Cells.Item(7, 13).Value = 36.48539666666667
$ws.Cells.Item(7, 14).Value = 109.45619
$ws.Cells.Item(7, 15).Value = 0.4260639713374229
$ws.Cells.Item(7, 16).Value = 0.4324607845540777
$ws.Cells.Item(7, 17).Value = 62.02445678719889
$ws.Cells.Item(7, 18).Value = 558.2201110847901
$ws.Cells.Item(7, 19).Value = 0.01676120619812524
$ws.Cells.Item(7, 20).Value = 0.01718271878513324
$ws.Cells.Item(8, 9).Value = 0.03933964692088724
$ws.Cells.Item(8, 10).Value = 0.0397324321622614
$ws.Cells.Item(8, 15).Value = 0.1743777127077069
$ws.Cells.Item(8, 16).Value = 0.1769957741547643
$ws.Cells.Item(8, 19).Value = 0.0068599576487931
$ws.Cells.Item(8, 20).Value = 0.007032472589611113
$ws.Cells.Item(9, 9).Value = 0.03933964692088724
$ws.Cells.Item(9, 10).Value = 0.0397324321622614
$ws.Cells.Item(9, 13).Value = 9.680823666666667
$ws.Cells.Item(9, 14).Value = 29.042471
$ws.Cells.Item(9, 15).Value = 0.1130493445068016
$ws.Cells.Item(9, 16).Value = 0.1147466378470605
$ws.Cells.Item(9, 17).Value = 16.45720984380122
$ws.Cells.Item(9, 18).Value = 148.114888594211
$ws.Cells.Item(9, 19).Value = 0.004447321297535319
$ws.Cells.Item(9, 20).Value = 0.00455916300410591
$ws.Cells.Item(10, 9).Value = 0.03933964692088724
$ws.Cells.Item(10, 10).Value = 0.0397324321622614
$ws.Cells.Item(10, 13).Value = 3.79999
$ws.Cells.Item(10, 14).Value = 7.59998
$ws.Cells.Item(10, 15).Value = 0.04437498227672168
$ws.Cells.Item(10, 16).Value = 0.0300274777826206
$ws.Cells.Item(10, 17).Value = 6.459908266863334
$ws.Cells.Item(10, 18).Value = 38.75944960118
$ws.Cells.Item(10, 19).Value = 0.00174569613488686
$ws.Cells.Item(10, 20).Value = 0.001193064724001784
$ws.Cells.Item(11, 9).Value = 0.03933964692088724
$ws.Cells.Item(11, 10).Value = 0.0397324321622614
$ws.Cells.Item(11, 13).Value = 20.734808
$ws.Cells.Item(11, 14).Value = 62.204424
$ws.Cells.Item(11, 15).Value = 0.242133989171347
$ws.Cells.Item(11, 16).Value = 0.245769325661477
$ws.Cells.Item(11, 17).Value = 35.24876581544267
$ws.Cells.Item(11, 18).Value = 317.238892338984
$ws.Cells.Item(11, 19).Value = 0.009525465641546726
$ws.Cells.Item(11, 20).Value = 0.009765013059409364
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.143548
$ws.Cells.Item(12, 8).Value = 0.430644
$ws.Cells.Item(12, 9).Value = 0.003321878215571232
$ws.Cells.Item(12, 10).Value = 0.003355045385051493
$ws.Cells.Item(12, 13).Value = 36.48539666666667
$ws.Cells.Item(12, 14).Value = 109.45619
$ws.Cells.Item(12, 15).Value = 0.4260639713374229
$ws.Cells.Item(12, 16).Value = 0.4324607845540777
$ws.Cells.Item(12, 17).Value = 5.237405720706667
$ws.Cells.Item(12, 18).Value = 47.13665148636001
$ws.Cells.Item(12, 19).Value = 0.001415332624825551
$ws.Cells.Item(12, 20).Value = 0.001450925559433906
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.143548
$ws.Cells.Item(13, 8).Value = 0.430644
$ws.Cells.Item(13, 9).Value = 0.003321878215571232
$ws.Cells.Item(13, 10).Value = 0.003355045385051493
$ws.Cells.Item(13, 15).Value = 0.1743777127077069
$ws.Cells.Item(13, 16).Value = 0.1769957741547643
$ws.Cells.Item(13, 17).Value = 2.143543907813334
$ws.Cells.Item(13, 18).Value = 19.29189517032
$ws.Cells.Item(13, 19).Value = 0.0005792615251248702
$ws.Cells.Item(13, 20).Value = 0.0005938288552515584
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.143548
$ws.Cells.Item(14, 8).Value = 0.430644
$ws.Cells.Item(14, 9).Value = 0.003321878215571232
$ws.Cells.Item(14, 10).Value = 0.003355045385051493
$ws.Cells.Item(14, 13).Value = 9.680823666666667
$ws.Cells.Item(14, 14).Value = 29.042471
$ws.Cells.Item(14, 15).Value = 0.1130493445068016
$ws.Cells.Item(14, 16).Value = 0.1147466378470605
$ws.Cells.Item(14, 17).Value = 1.389662875702667
$ws.Cells.Item(14, 18).Value = 12.506965881324
$ws.Cells.Item(14, 19).Value = 0.0003755361548017516
$ws.Cells.Item(14, 20).Value = 0.0003849801777589554
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.143548
$ws.Cells.Item(15, 8).Value = 0.430644
$ws.Cells.Item(15, 9).Value = 0.003321878215571232
$ws.Cells.Item(15, 10).Value = 0.003355045385051493
$ws.Cells.Item(15, 13).Value = 3.79999
$ws.Cells.Item(15, 14).Value = 7.59998
$ws.Cells.Item(15, 15).Value = 0.04437498227672168
$ws.Cells.Item(15, 16).Value = 0.0300274777826206
$ws.Cells.Item(15, 17).Value = 0.5454809645200001
$ws.Cells.Item(15, 18).Value = 3.27288578712
$ws.Cells.Item(15, 19).Value = 0.0001474082869414012
$ws.Cells.Item(15, 20).Value = 0.0001007435507593175
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.143548
$ws.Cells.Item(16, 8).Value = 0.430644
$ws.Cells.Item(16, 9).Value = 0.003321878215571232
$ws.Cells.Item(16, 10).Value = 0.003355045385051493
$ws.Cells.Item(16, 13).Value = 20.734808
$ws.Cells.Item(16, 14).Value = 62.204424
$ws.Cells.Item(16, 15).Value = 0.242133989171347
$ws.Cells.Item(16, 16).Value = 0.245769325661477
$ws.Cells.Item(16, 17).Value = 2.976440218784
$ws.Cells.Item(16, 18).Value = 26.787961969056
$ws.Cells.Item(16, 19).Value = 0.0008043396238776581
$ws.Cells.Item(16, 20).Value = 0.0008245672418477559
$ws.Cells.Item(17, 7).Value = 1.281577
$ws.Cells.Item(17, 8).Value = 2.563154
$ws.Cells.Item(17, 9).Value = 0.02965727643629401
$ws.Cells.Item(17, 10).Value = 0.01996892560647838
$ws.Cells.Item(17, 13).Value = 36.48539666666667
$ws.Cells.Item(17, 14).Value = 109.45619
$ws.Cells.Item(17, 15).Value = 0.4260639713374229
$ws.Cells.Item(17, 16).Value = 0.4324607845540777
$ws.Cells.Item(17, 17).Value = 46.75884520387667
$ws.Cells.Item(17, 18).Value = 280.55307122326
$ws.Cells.Item(17, 19).Value = 0.0126358969774992
$ws.Cells.Item(17, 20).Value = 0.008635777234479652
$ws.Cells.Item(18, 7).Value = 1.281577
$ws.Cells.Item(18, 8).Value = 2.563154
$ws.Cells.Item(18, 9).Value = 0.02965727643629401
$ws.Cells.Item(18, 10).Value = 0.01996892560647838
$ws.Cells.Item(18, 15).Value = 0.1743777127077069
$ws.Cells.Item(18, 16).Value = 0.1769957741547643
$ws.Cells.Item(18, 17).Value = 19.13726816635333
$ws.Cells.Item(18, 18).Value = 114.82360899812
$ws.Cells.Item(18, 19).Value = 0.005171568030101121
$ws.Cells.Item(18, 20).Value = 0.003534415446757537
$ws.Cells.Item(19, 7).Value = 1.281577
$ws.Cells.Item(19, 8).Value = 2.563154
$ws.Cells.Item(19, 9).Value = 0.02965727643629401
$ws.Cells.Item(19, 10).Value = 0.01996892560647838
$ws.Cells.Item(19, 13).Value = 9.680823666666667
$ws.Cells.Item(19, 14).Value = 29.042471
$ws.Cells.Item(19, 15).Value = 0.1130493445068016
$ws.Cells.Item(19, 16).Value = 0.1147466378470605
$ws.Cells.Item(19, 17).Value = 12.40672095225567
$ws.Cells.Item(19, 18).Value = 74.440325713534
$ws.Cells.Item(19, 19).Value = 0.003352735660980051
$ws.Cells.Item(19, 20).Value = 0.002291367074761468
$ws.Cells.Item(20, 7).Value = 1.281577
$ws.Cells.Item(20, 8).Value = 2.563154
$ws.Cells.Item(20, 9).Value = 0.02965727643629401
$ws.Cells.Item(20, 10).Value = 0.01996892560647838
$ws.Cells.Item(20, 13).Value = 3.79999
$ws.Cells.Item(20, 14).Value = 7.59998
$ws.Cells.Item(20, 15).Value = 0.04437498227672168
$ws.Cells.Item(20, 16).Value = 0.0300274777826206
$ws.Cells.Item(20, 17).Value = 4.86997978423
$ws.Cells.Item(20, 18).Value = 19.47991913692
$ws.Cells.Item(20, 19).Value = 0.001316041116236382
$ws.Cells.Item(20, 20).Value = 0.0005996164699913332
$ws.Cells.Item(21, 7).Value = 1.281577
$ws.Cells.Item(21, 8).Value = 2.563154
$ws.Cells.Item(21, 9).Value = 0.02965727643629401
$ws.Cells.Item(21, 10).Value = 0.01996892560647838
$ws.Cells.Item(21, 13).Value = 20.734808
$ws.Cells.Item(21, 14).Value = 62.204424
$ws.Cells.Item(21, 15).Value = 0.242133989171347
$ws.Cells.Item(21, 16).Value = 0.245769325661477
$ws.Cells.Item(21, 17).Value = 26.573253032216
$ws.Cells.Item(21, 18).Value = 159.439518193296
$ws.Cells.Item(21, 19).Value = 0.007181034651477259
$ws.Cells.Item(21, 20).Value = 0.004907749380488391
$ws.Cells.Item(22, 7).Value = 0.07451633333333334
$ws.Cells.Item(22, 8).Value = 0.223549
$ws.Cells.Item(22, 9).Value = 0.001724400091984872
$ws.Cells.Item(22, 10).Value = 0.001741617300561197
$ws.Cells.Item(22, 13).Value = 36.48539666666667
$ws.Cells.Item(22, 14).Value = 109.45619
$ws.Cells.Item(22, 15).Value = 0.4260639713374229
$ws.Cells.Item(22, 16).Value = 0.4324607845540777
$ws.Cells.Item(22, 17).Value = 2.718757979812223
$ws.Cells.Item(22, 18).Value = 24.46882181831
$ws.Cells.Item(22, 19).Value = 0.000734704751365692
$ws.Cells.Item(22, 20).Value = 0.0007531811841936504
$ws.Cells.Item(23, 7).Value = 0.07451633333333334
$ws.Cells.Item(23, 8).Value = 0.223549
$ws.Cells.Item(23, 9).Value = 0.001724400091984872
$ws.Cells.Item(23, 10).Value = 0.001741617300561197
$ws.Cells.Item(23, 15).Value = 0.1743777127077069
$ws.Cells.Item(23, 16).Value = 0.1769957741547643
$ws.Cells.Item(23, 17).Value = 1.112722102357778
$ws.Cells.Item(23, 18).Value = 10.01449892122
$ws.Cells.Item(23, 19).Value = 0.0003006969438332813
$ws.Cells.Item(23, 20).Value = 0.00030825890239416
$ws.Cells.Item(24, 7).Value = 0.07451633333333334
$ws.Cells.Item(24, 8).Value = 0.223549
$ws.Cells.Item(24, 9).Value = 0.001724400091984872
$ws.Cells.Item(24, 10).Value = 0.001741617300561197
$ws.Cells.Item(24, 13).Value = 9.680823666666667
$ws.Cells.Item(24, 14).Value = 29.042471
$ws.Cells.Item(24, 15).Value = 0.1130493445068016
$ws.Cells.Item(24, 16).Value = 0.1147466378470605
$ws.Cells.Item(24, 17).Value = 0.7213794832865557
$ws.Cells.Item(24, 18).Value = 6.492415349579
$ws.Cells.Item(24, 19).Value = 0.0001949423000663582
$ws.Cells.Item(24, 20).Value = 0.0001998447296556709
$ws.Cells.Item(25, 7).Value = 0.07451633333333334
$ws.Cells.Item(25, 8).Value = 0.223549
$ws.Cells.Item(25, 9).Value = 0.001724400091984872
$ws.Cells.Item(25, 10).Value = 0.001741617300561197
$ws.Cells.Item(25, 13).Value = 3.79999
$ws.Cells.Item(25, 14).Value = 7.59998
$ws.Cells.Item(25, 15).Value = 0.04437498227672168
$ws.Cells.Item(25, 16).Value = 0.0300274777826206
$ws.Cells.Item(25, 17).Value = 0.2831613215033333
$ws.Cells.Item(25, 18).Value = 1.69896792902
$ws.Cells.Item(25, 19).Value = 0.00007652022351980594
$ws.Cells.Item(25, 20).Value = 0.00005229637479842902
$ws.Cells.Item(26, 7).Value = 0.07451633333333334
$ws.Cells.Item(26, 8).Value = 0.223549
$ws.Cells.Item(26, 9).Value = 0.001724400091984872
$ws.Cells.Item(26, 10).Value = 0.001741617300561197
$ws.Cells.Item(26, 13).Value = 20.734808
$ws.Cells.Item(26, 14).Value = 62.204424
$ws.Cells.Item(26, 15).Value = 0.242133989171347
$ws.Cells.Item(26, 16).Value = 0.245769325661477
$ws.Cells.Item(26, 17).Value = 1.545081864530667
$ws.Cells.Item(26, 18).Value = 13.905736780776
$ws.Cells.Item(26, 19).Value = 0.0004175358731997349
$ws.Cells.Item(26, 20).Value = 0.0004280361095192873
